# Tp Nro 5 terminado
# Update the ranking table (Score / Documento / Contenido) on the active sheet:
#  - row 2 gets a new score, a new source document name, and new content
#  - row 3 only gets an updated score
#  - row 4 gets a new score, a new source document name, and new content
#
# The Score column stores its numbers as text (e.g. "0.089"), so we force the
# NumberFormat to Text before writing, then clear the format again so the
# cell keeps the workbook's original (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0.102"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = 'Noticia N° 09.txt'
$ws.Range("D2").Value = 'título: aporte dólar soja, banco central compró solo usd millones mercado resumen: entidad monetaria encadenó ruedas operativas compras netas, margen achicó primer día agregado tipo cambio especial exportaciones sojeras contenido: jun, sesión usd millones segmento contado spot, liquidaciones programa incremento exportador economías regionales, dólar, aportaron usd millones, volumen aportado ventas soja derivados, gozaron estándar cambiario semana pasada. bcra concluyó intervención cambiaria día saldo neto comprador solo millones dólares. operaciones dólar agro registran tercera etapa ingresos usd millones abril parte. asimismo, empezó balance bcra intervención cambiaria negativo usd millones, mientras junio sostiene im saldo comprador millones dólares. bcra aceleró resguardo escasas reservas, medida afecta provincias municipios, momentos ministerio economía planea nuevo canje voluntario bonos pesos objetivo despejar vencimientos. ministro economía, sergio massa, viajará washington cerca junio procura cerrar readecuación acuerdo fondo monetario internacional fmi, permitiría argentina obtener adelanto desembolsos nuevas metas cumplir, luego impacto sequía exportaciones sector agro, “las negociaciones argentina funcionarios técnicos fondo avanzan hace casi dos meses través reuniones virtuales, definiendo modificar viejo acuerdo, virtualmente suspendido luego cumplieran metas reservas fiscales primer trimestre ″, comentaron expertos research traders. objetivo economía fondo adelante, menos, parte desembolsos comprometidos fin año usd millones, ayudaría reforzar reservas año exportaciones sector agro caerían cerca usd millones, acuerdo proyecciones bolsa comercio rosario bcr. tales desembolsos dudas cuánto podría utilizar intervenir mercado cambiario. fondo aceptaría bcra intervenga eventuales situaciones stress. todavía resuelto monto fmi consentiría intervenciones. reservas internacionales brutas banco central crecieron semana pasada usd millones finalizaron millones dólares. informe anker latinoamérica subrayó cuanto reservas líquidas bcra, “su disponibilidad hoy dada gran medida encajes cuentas bancarias monedas depositados bcra -usd millones-. estabilidad depósitos moneda extranjera crucial sostener capacidad intervención bcra”. luego diversas gestiones ministro massa incentivar pago importaciones divisa china, expandió uso yuanes comercio exterior. enero mayo operaciones equivalente usd millones. además, dos empresas fabricantes electrónica confirmaron pagarán compromisos deuda total usd millones moneda, alivio arcas bcra superará usd millones primer trimestre seguir leyendo: urls imagenes:'

# Row 3
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0.069"
$ws.Range("B3").ClearFormats()

# Row 4
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "0.054"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = 'Noticia N° 06.txt'
$ws.Range("D4").Value = 'título: inflación mayo argentina superó venezuela, según estimaciones privadas resumen: abril registrado resultado; variado precios último año ambos países espera diciembre contenido: jun, martín kanenguiser inflación argentina mayo superó venezuela, vez más, según estimaciones privadas ambos países. dato observatorio venezolano finanzas ovf, siguen analistas debido falta rigor estadísticas oficiales, marcó suba precios mes pasado, desaceleración últimos meses. argentina, según mayoría estimaciones, dato mes pasado rondó %, mientras gobierno reza ubique levemente debajo cifra. abril, argentina superado venezuela: registró %, frente informado ovf informado banco central venezuela, daba conocer información octubre año pasado. últimos meses mayo-mayo, según ovf, inflación país gobernado nicolás maduro llegó %, baja respecto registrado abril. además, enero suba precios acumulada %. tendencias inflación últimos meses divergentes: marzo ovf registró venezuela %, abril mayo %; tanto, indec registró argentina marzo, abril estima mes quedará torno interanual; conocerá dato organismo lidera marco lavagna. según ovf, “en medio economía evidentes signos contracción, tasa inflación da tregua. así, mayo tasa inflación mensual alcanzó %, acumulada anualizada %. respecto abril, tasa mensual inflación triplicó aunque tasa doce meses desaceleró”. “este comportamiento inflación ocurriendo entorno signos significativa caída nivel actividad indiscutibles, ostensibles menores ventas comercio producción industrial contrajo primer trimestre ″, indicó. además, “la debilidad demanda agregada notoria debido salarios pensiones sector público pronunciadamente rezagados respecto inflación ejecución gasto parte gobierno baja”. “al comparar cifra inflación mayo respecto aumento tipo cambio mismo mes, clara estrecha relación ambas variables, obstante, mencionó, debilidad demanda, explicada política compresión salarial venido aplicando gobierno”, subrayó ovf. cuanto principales componentes conforman índice nacional precios consumidor, destacaron “los aumentos experimentados rubros esparcimiento %, vestido calzado %, equipamiento hogar alquiler vivienda %. alimentos alzas modestas, mayo incrementaron %”. “todo ello pone manifiesto, política económica aplicada sido ineficaz contener alza precios, aún retracción inducida demanda agregada”, concluyó organismo independiente régimen autoritario maduro. abril, banco central venezuela informado aumento cuatro meses, luego difundir datos medio año. tanto, mayoría relevamientos precios consumidor consultoras argentinas anticipan ipc torno %. c&t indicó relevamiento precios minoristas región gba “presentó alza mensual %, superando largamente variación abril mayo año pasado. así, variación doce meses trepó %, mayor agosto ″. “el rubro mayor incremento vivienda %, reflejando subas gas electricidad principalmente”, aclaró. bienes servicios varios “ocuparon segundo lugar, alza %, explicada cigarrillos artículos tocador”. tanto, “el comportamiento esparcimiento fuertemente influido alza dólares financieros fines abril, vio reflejado turismo productos electrónicos”. vez, “alimentos bebidas creció mes. arrancó mes gran impulso luego moderando. verduras, lácteos derivados harina destacaron, igual alimentos consumidos hogar llevar”. “en salud destacó incremento medicamentos, sumó nuevo ajuste prepagas”, indicó c&t. “en equipamiento mantenimiento hogar, artefactos hogar rol preponderante mano alza dólares financieros pesar liquidaciones mes”, concluyó. parte, ecogo informó si bien última semana mes pasado registró fuerte desaceleración lugar variación alimentos respecto semana anterior, términos generales inflación sido ciento. particular, aumento precios alimentos sido ciento. “si consideramos además evolución alimentos consumidos hogar %, inflación alimentos alcanzaría %”, aclaró consultora dirige marina dal poggeto. lado, lcg detalló “el índice alimentos bebidas presentó inflación mensual promedio últimas semanas punta punta mismo período”. mes sumarán aumentos precios servicios transporte, prepagas, colegios privados combustibles, valores regulados. equipo económico trata controlar cuestión cambiaria inflación acelere todavía motivo resiste pedido fmi acelerar devaluación tipo cambio oficial. seguir leyendo: urls imagenes:'
